# The workbook tracks weekly price observations for "Frutilla" (strawberry)
# at Feria Lagunitas de Puerto Montt. This edit adds two new weekly price
# records (dated 2022-01-11 / serial 44572) at the top of the historical
# data block (rows 119-120), pushing all the existing historical rows
# down by two rows (old row 119 -> new row 121, ..., old row 177 -> new
# row 179). The sheet's dimension grows from A1:T177 to A1:T179.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the old row 119, shifting the
# remaining historical rows (old 119..177) down to 121..179.
$ws.Rows("119:120").Insert()

# --- New row 119 ---------------------------------------------------------
$ws.Range("A119").Value = 4
$ws.Range("B119").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C119").Value = "Los Lagos"
$ws.Range("D119").Value = 44572
$ws.Range("E119").Value = 10
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100101
$ws.Range("H119").Value = "Berries"
$ws.Range("I119").Value = 100112025
$ws.Range("J119").Value = "Frutilla"
$ws.Range("K119").Value = "Sin especificar"
$ws.Range("L119").Value = "Primera"
$ws.Range("M119").Value = 800
$ws.Range("N119").Value = 8500
$ws.Range("O119").Value = 9000
$ws.Range("P119").Value = 8750
$ws.Range("Q119").Value = "`$/caja 7 kilos"
$ws.Range("R119").Value = "Región de La Araucanía"
$ws.Range("S119").Value = 1250
$ws.Range("T119").Value = 7

# --- New row 120 ---------------------------------------------------------
$ws.Range("A120").Value = 4
$ws.Range("B120").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C120").Value = "Los Lagos"
$ws.Range("D120").Value = 44572
$ws.Range("E120").Value = 10
$ws.Range("F120").Value = "Fruta"
$ws.Range("G120").Value = 100101
$ws.Range("H120").Value = "Berries"
$ws.Range("I120").Value = 100112025
$ws.Range("J120").Value = "Frutilla"
$ws.Range("K120").Value = "Sin especificar"
$ws.Range("L120").Value = "Segunda"
$ws.Range("M120").Value = 300
$ws.Range("N120").Value = 7000
$ws.Range("O120").Value = 7000
$ws.Range("P120").Value = 7000
$ws.Range("Q120").Value = "`$/caja 7 kilos"
$ws.Range("R120").Value = "Región de La Araucanía"
$ws.Range("S120").Value = 1000
$ws.Range("T120").Value = 7

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
